# Applies the "Updated cryptos list" data refresh (Fri Jul  7 03:14:59 UTC 2023).
# Columns: A=rank index, B=Coin, C=Link, D=Price, E=Volume(1h). All of these are
# stored as plain text in the sheet, so numeric-looking Price values are written
# with a leading apostrophe to stop Excel from re-interpreting them as numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.093.64"
$ws.Range("E2").Value = "  -1.23%  "
$ws.Range("D3").Value = "1.854.65"
$ws.Range("E3").Value = "  -2.97%  "
$ws.Range("D4").Value = "'1.000"
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").Value = "'234.01"
$ws.Range("E5").Value = "  -2.07%  "
$ws.Range("D6").Value = "'1.000"
$ws.Range("E6").Value = "  +0.07%  "
$ws.Range("D7").Value = "'0.4689"
$ws.Range("E7").Value = "  -1.74%  "
$ws.Range("E8").Value = "  -0.85%  "
$ws.Range("D9").Value = "'0.06555"
$ws.Range("E9").Value = "  -2.12%  "
$ws.Range("D10").Value = "'20.20"
$ws.Range("E10").Value = "  +3.40%  "
$ws.Range("D11").Value = "'0.07793"
$ws.Range("D12").Value = "'97.12"
$ws.Range("E12").Value = "  -6.01%  "
$ws.Range("D13").Value = "1.852.57"
$ws.Range("E13").Value = "  -3.13%  "
$ws.Range("D14").Value = "'5.074"
$ws.Range("E14").Value = "  -2.20%  "
$ws.Range("D15").Value = "'0.6686"
$ws.Range("E15").Value = "  +0.23%  "
$ws.Range("D16").Value = "'283.64"
$ws.Range("E16").Value = "  +2.84%  "
$ws.Range("D17").Value = "30.112.35"
$ws.Range("E17").Value = "  -1.27%  "
$ws.Range("D18").Value = "'1.000"
$ws.Range("E18").Value = "  +0.04%  "
$ws.Range("D19").Value = "'12.59"
$ws.Range("E19").Value = "  -0.16%  "
$ws.Range("D20").Value = "'5.397"
$ws.Range("E20").Value = "  +0.38%  "
$ws.Range("B21").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C21").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D21").Value = "2.098.82"
$ws.Range("E21").Value = "  -2.67%  "
$ws.Range("B22").Value = "ShibaInu"
$ws.Range("C22").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D22").Value = "'0.000007232"
$ws.Range("E22").Value = "  -3.29%  "
$ws.Range("D23").Value = "'1.000"
$ws.Range("E23").Value = "  +0.08%  "
$ws.Range("D24").Value = "'6.132"
$ws.Range("E24").Value = "  -2.42%  "
$ws.Range("D25").Value = "'167.59"
$ws.Range("E25").Value = "  +0.65%  "
$ws.Range("D26").Value = "'9.286"
$ws.Range("E26").Value = "  -0.45%  "
$ws.Range("D27").Value = "'18.99"
$ws.Range("E27").Value = "  -1.09%  "
$ws.Range("D28").Value = "'1.925"
$ws.Range("E28").Value = "  -7.20%  "
$ws.Range("D29").Value = "'1.336"
$ws.Range("E29").Value = "  -3.08%  "
$ws.Range("D30").Value = "'0.09673"
$ws.Range("E30").Value = "  -3.07%  "
$ws.Range("D31").Value = "'4.378"
$ws.Range("E31").Value = "  -5.01%  "
$ws.Range("D32").Value = "'1.468"
$ws.Range("E32").Value = "  -2.60%  "
$ws.Range("D33").Value = "'4.094"
$ws.Range("E33").Value = "  -3.62%  "
$ws.Range("D34").Value = "'0.04658"
$ws.Range("D35").Value = "'0.6982"
$ws.Range("E35").Value = "  -3.85%  "
$ws.Range("D36").Value = "'1.084"
$ws.Range("E36").Value = "  -2.52%  "
$ws.Range("D37").Value = "'0.9997"
$ws.Range("E37").Value = "  +0.10%  "
$ws.Range("D38").Value = "'2.699"
$ws.Range("E38").Value = "  -0.35%  "
$ws.Range("E39").Value = "  -2.71%  "
$ws.Range("D40").Value = "'6.299"
$ws.Range("E40").Value = "  -0.46%  "
$ws.Range("D41").Value = "'2.503"
$ws.Range("E41").Value = "  -3.99%  "
$ws.Range("D42").Value = "'71.53"
$ws.Range("E42").Value = "  -3.72%  "
$ws.Range("D43").Value = "'0.8600"
$ws.Range("E43").Value = "  -0.03%  "
$ws.Range("D44").Value = "'1.934"
$ws.Range("E44").Value = "  -1.23%  "
$ws.Range("D45").Value = "'104.17"
$ws.Range("E45").Value = "  -1.57%  "
$ws.Range("D46").Value = "'1.000"
$ws.Range("E46").Value = "  +0.09%  "
$ws.Range("D47").Value = "'0.4157"
$ws.Range("E47").Value = "  -2.34%  "
$ws.Range("D48").Value = "1.015.71"
$ws.Range("E48").Value = "  +6.42%  "
$ws.Range("D49").Value = "'7.247"
$ws.Range("E49").Value = "  -2.15%  "
$ws.Range("D50").Value = "'9.216"
$ws.Range("E50").Value = "  +4.82%  "
$ws.Range("D51").Value = "'33.76"
$ws.Range("E51").Value = "  -2.39%  "
